$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.465.31'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.883.38'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '243.75'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4716'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2882'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06465'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '22.17'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.77%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07771'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.887.17'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '95.53'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.7223'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.172'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '280.13'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.00%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.468.24'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.03'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.55%  '
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007452'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.134.67'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.99%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.251'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.72%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.277'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '164.00'
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.050'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.43%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.79'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.886'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09630'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.25%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.469'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.50%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.256'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.132'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.38%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04846'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.123'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6911'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.710'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01881'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.62%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.820'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.44%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '74.53'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.72%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.194'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.81%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.958'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.4257'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.54%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8262'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.89%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '101.07'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.628'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.08%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '6.947'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.21%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '35.15'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '903.81'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.26%  '
$ws.Range('E51').Value = '  +1.57%  '
